$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-51: refresh Price (D) and Volume(1h) (E) columns.
# Price strings look numeric (e.g. "290.97") so force text via
# a temporary "@" number format, then restore the default style
# so no stray formatting is left behind (matches source: plain
# inline/shared strings with no explicit cell style).

$d = $ws.Range("D2")
$d.NumberFormat = "@"
$d.Value = "22.402.16"
$d.Style = "Normal"
$ws.Range("E2").Value = "  -4.23%  "

$d = $ws.Range("D3")
$d.NumberFormat = "@"
$d.Value = "1.567.82"
$d.Style = "Normal"
$ws.Range("E3").Value = "  -3.95%  "

$ws.Range("E4").Value = "  +0.31%  "

$d = $ws.Range("D5")
$d.NumberFormat = "@"
$d.Value = "1.002"
$d.Style = "Normal"
$ws.Range("E5").Value = "  +0.24%  "

$d = $ws.Range("D6")
$d.NumberFormat = "@"
$d.Value = "290.97"
$d.Style = "Normal"
$ws.Range("E6").Value = "  -2.47%  "

$d = $ws.Range("D7")
$d.NumberFormat = "@"
$d.Value = "0.3691"
$d.Style = "Normal"
$ws.Range("E7").Value = "  -2.22%  "

$d = $ws.Range("D8")
$d.NumberFormat = "@"
$d.Value = "49.36"
$d.Style = "Normal"
$ws.Range("E8").Value = "  -1.58%  "

$d = $ws.Range("D9")
$d.NumberFormat = "@"
$d.Value = "0.3392"
$d.Style = "Normal"
$ws.Range("E9").Value = "  -2.50%  "

$d = $ws.Range("D10")
$d.NumberFormat = "@"
$d.Value = "1.167"
$d.Style = "Normal"
$ws.Range("E10").Value = "  -3.47%  "

$d = $ws.Range("D11")
$d.NumberFormat = "@"
$d.Value = "0.07652"
$d.Style = "Normal"
$ws.Range("E11").Value = "  -4.94%  "

$d = $ws.Range("D12")
$d.NumberFormat = "@"
$d.Value = "1.003"
$d.Style = "Normal"
$ws.Range("E12").Value = "  +0.45%  "

$d = $ws.Range("D13")
$d.NumberFormat = "@"
$d.Value = "21.37"
$d.Style = "Normal"
$ws.Range("E13").Value = "  -2.67%  "

$d = $ws.Range("D14")
$d.NumberFormat = "@"
$d.Value = "6.056"
$d.Style = "Normal"
$ws.Range("E14").Value = "  -3.74%  "

$d = $ws.Range("D15")
$d.NumberFormat = "@"
$d.Value = "6.914"
$d.Style = "Normal"
$ws.Range("E15").Value = "  -4.27%  "

$d = $ws.Range("D16")
$d.NumberFormat = "@"
$d.Value = "1.572.83"
$d.Style = "Normal"
$ws.Range("E16").Value = "  -3.58%  "

$d = $ws.Range("D17")
$d.NumberFormat = "@"
$d.Value = "0.00001128"
$d.Style = "Normal"
$ws.Range("E17").Value = "  -5.49%  "

$d = $ws.Range("D18")
$d.NumberFormat = "@"
$d.Value = "90.02"
$d.Style = "Normal"
$ws.Range("E18").Value = "  -4.95%  "

$d = $ws.Range("D19")
$d.NumberFormat = "@"
$d.Value = "0.06722"
$d.Style = "Normal"
$ws.Range("E19").Value = "  -3.21%  "

$d = $ws.Range("D20")
$d.NumberFormat = "@"
$d.Value = "1.002"
$d.Style = "Normal"
$ws.Range("E20").Value = "  +0.08%  "

$d = $ws.Range("D21")
$d.NumberFormat = "@"
$d.Value = "6.261"
$d.Style = "Normal"
$ws.Range("E21").Value = "  -5.31%  "

$d = $ws.Range("D22")
$d.NumberFormat = "@"
$d.Value = "16.54"
$d.Style = "Normal"
$ws.Range("E22").Value = "  -4.25%  "

$d = $ws.Range("D23")
$d.NumberFormat = "@"
$d.Value = "0.5311"
$d.Style = "Normal"
$ws.Range("E23").Value = "  -6.72%  "

$d = $ws.Range("D24")
$d.NumberFormat = "@"
$d.Value = "12.01"
$d.Style = "Normal"
$ws.Range("E24").Value = "  -2.96%  "

$d = $ws.Range("D25")
$d.NumberFormat = "@"
$d.Value = "22.377.68"
$d.Style = "Normal"
$ws.Range("E25").Value = "  -4.33%  "

$d = $ws.Range("D26")
$d.NumberFormat = "@"
$d.Value = "2.356"
$d.Style = "Normal"
$ws.Range("E26").Value = "  -2.54%  "

$d = $ws.Range("D27")
$d.NumberFormat = "@"
$d.Value = "2.850"
$d.Style = "Normal"
$ws.Range("E27").Value = "  -3.52%  "

$d = $ws.Range("D28")
$d.NumberFormat = "@"
$d.Value = "20.06"
$d.Style = "Normal"
$ws.Range("E28").Value = "  -4.08%  "

$d = $ws.Range("D29")
$d.NumberFormat = "@"
$d.Value = "146.40"
$d.Style = "Normal"
$ws.Range("E29").Value = "  -1.92%  "

$d = $ws.Range("D30")
$d.NumberFormat = "@"
$d.Value = "4.981"
$d.Style = "Normal"
$ws.Range("E30").Value = "  -3.83%  "

$d = $ws.Range("D31")
$d.NumberFormat = "@"
$d.Value = "125.68"
$d.Style = "Normal"
$ws.Range("E31").Value = "  -4.03%  "

$d = $ws.Range("D32")
$d.NumberFormat = "@"
$d.Value = "1.745.48"
$d.Style = "Normal"
$ws.Range("E32").Value = "  -3.36%  "

$d = $ws.Range("D33")
$d.NumberFormat = "@"
$d.Value = "1.019"
$d.Style = "Normal"
$ws.Range("E33").Value = "  +3.64%  "

$d = $ws.Range("D34")
$d.NumberFormat = "@"
$d.Value = "6.194"
$d.Style = "Normal"
$ws.Range("E34").Value = "  -8.09%  "

$d = $ws.Range("D35")
$d.NumberFormat = "@"
$d.Value = "2.019"
$d.Style = "Normal"
$ws.Range("E35").Value = "  -4.55%  "

$d = $ws.Range("D36")
$d.NumberFormat = "@"
$d.Value = "10.08"
$d.Style = "Normal"
$ws.Range("E36").Value = "  -9.45%  "

$d = $ws.Range("D37")
$d.NumberFormat = "@"
$d.Value = "0.08536"
$d.Style = "Normal"
$ws.Range("E37").Value = "  -2.30%  "

$d = $ws.Range("D38")
$d.NumberFormat = "@"
$d.Value = "0.02539"
$d.Style = "Normal"
$ws.Range("E38").Value = "  -4.48%  "

$d = $ws.Range("D39")
$d.NumberFormat = "@"
$d.Value = "0.2317"
$d.Style = "Normal"
$ws.Range("E39").Value = "  -4.00%  "

$d = $ws.Range("D40")
$d.NumberFormat = "@"
$d.Value = "5.518"
$d.Style = "Normal"
$ws.Range("E40").Value = "  -5.34%  "

$d = $ws.Range("D41")
$d.NumberFormat = "@"
$d.Value = "0.06469"
$d.Style = "Normal"
$ws.Range("E41").Value = "  -4.50%  "

$d = $ws.Range("D42")
$d.NumberFormat = "@"
$d.Value = "1.277"
$d.Style = "Normal"
$ws.Range("E42").Value = "  -0.94%  "

$d = $ws.Range("D43")
$d.NumberFormat = "@"
$d.Value = "11.68"
$d.Style = "Normal"
$ws.Range("E43").Value = "  -8.37%  "

$d = $ws.Range("D44")
$d.NumberFormat = "@"
$d.Value = "0.6339"
$d.Style = "Normal"
$ws.Range("E44").Value = "  -6.71%  "

$d = $ws.Range("D45")
$d.NumberFormat = "@"
$d.Value = "14.24"
$d.Style = "Normal"
$ws.Range("E45").Value = "  -7.32%  "

$d = $ws.Range("D46")
$d.NumberFormat = "@"
$d.Value = "1.001"
$d.Style = "Normal"
$ws.Range("E46").Value = "  +0.19%  "

$d = $ws.Range("D47")
$d.NumberFormat = "@"
$d.Value = "0.5981"
$d.Style = "Normal"
$ws.Range("E47").Value = "  -5.38%  "

$d = $ws.Range("D48")
$d.NumberFormat = "@"
$d.Value = "3.762"
$d.Style = "Normal"
$ws.Range("E48").Value = "  -3.28%  "

$d = $ws.Range("D49")
$d.NumberFormat = "@"
$d.Value = "2.108"
$d.Style = "Normal"
$ws.Range("E49").Value = "  -5.39%  "

$d = $ws.Range("D50")
$d.NumberFormat = "@"
$d.Value = "1.260"
$d.Style = "Normal"
$ws.Range("E50").Value = "  +3.37%  "

$d = $ws.Range("D51")
$d.NumberFormat = "@"
$d.Value = "125.16"
$d.Style = "Normal"
$ws.Range("E51").Value = "  -1.01%  "
